# Add new EV region "UK00" to the Capacity sheet, mirroring the existing
# NOS0 block (rows 10-11) with new capacity figures for 2030/2040.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the two new rows by copying the last existing data row (11, NOS0 /
# 2040) down into rows 12 and 13. This carries over the column formatting
# (Scenario/Year/coefficient styles) exactly as the existing rows have it.
$ws.Rows("11").Copy()
$ws.Rows("12").Insert(-4121)
$ws.Rows("11").Copy()
$ws.Rows("13").Insert(-4121)

# Row 12: UK00, Distributed Energy, 2030, 200000
$ws.Range("A12").Value = "UK00"
$ws.Range("C12").Value = 2030
$ws.Range("D12").Value = 200000

# Row 13: UK00, Distributed Energy, 2040, 400000
$ws.Range("A13").Value = "UK00"
$ws.Range("C13").Value = 2040
$ws.Range("D13").Value = 400000

# Restore the active selection to where the user ended up after entering
# the new rows.
$ws.Range("A16").Select()
